$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates in column D between rows 2-3 and rows 6-7.
# D2, D3: 44574 -> 44559
# D6, D7: 44559 -> 44574
$ws.Range("D2").Value = 44559
$ws.Range("D3").Value = 44559
$ws.Range("D6").Value = 44574
$ws.Range("D7").Value = 44574
